$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamps = @(
    "Sat Oct 12 23:29:42 EDT 2024",
    "Sat Oct 12 23:29:53 EDT 2024",
    "Sat Oct 12 23:30:03 EDT 2024",
    "Sat Oct 12 23:30:13 EDT 2024",
    "Sat Oct 12 23:30:25 EDT 2024",
    "Sat Oct 12 23:30:35 EDT 2024",
    "Sat Oct 12 23:30:45 EDT 2024",
    "Sat Oct 12 23:30:56 EDT 2024",
    "Sat Oct 12 23:31:08 EDT 2024",
    "Sat Oct 12 23:31:19 EDT 2024",
    "Sat Oct 12 23:31:30 EDT 2024",
    "Sat Oct 12 23:31:41 EDT 2024",
    "Sat Oct 12 23:31:51 EDT 2024",
    "Sat Oct 12 23:32:01 EDT 2024",
    "Sat Oct 12 23:32:12 EDT 2024",
    "Sat Oct 12 23:32:23 EDT 2024",
    "Sat Oct 12 23:32:33 EDT 2024",
    "Sat Oct 12 23:32:44 EDT 2024",
    "Sat Oct 12 23:32:54 EDT 2024",
    "Sat Oct 12 23:33:05 EDT 2024",
    "Sat Oct 12 23:33:15 EDT 2024",
    "Sat Oct 12 23:33:25 EDT 2024",
    "Sat Oct 12 23:33:36 EDT 2024",
    "Sat Oct 12 23:33:46 EDT 2024",
    "Sat Oct 12 23:33:56 EDT 2024",
    "Sat Oct 12 23:34:07 EDT 2024",
    "Sat Oct 12 23:34:17 EDT 2024",
    "Sat Oct 12 23:34:28 EDT 2024",
    "Sat Oct 12 23:34:39 EDT 2024"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $newTimestamps[$i]
}
